$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2290909090909091
$ws.Range("C2").Value = 0.4690909090909091
$ws.Range("J2").Value = 0.007272727272727273
$ws.Range("P2").Value = 0.1781818181818182
$ws.Range("S2").Value = 0.1163636363636364
$ws.Range("B3").Value = 0.007518796992481203
$ws.Range("C3").Value = 0.02255639097744361
$ws.Range("J3").Value = 0.03007518796992481
$ws.Range("P3").Value = 0.706766917293233
$ws.Range("S3").Value = 0.2330827067669173
$ws.Range("J4").Value = 0.03333333333333333
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.2666666666666667
$ws.Range("B6").Value = 0.0390625
$ws.Range("D6").Value = 0.01953125
$ws.Range("F6").Value = 0.046875
$ws.Range("J6").Value = 0.24609375
$ws.Range("O6").Value = 0.02734375
$ws.Range("Q6").Value = 0.11328125
$ws.Range("R6").Value = 0.09375
$ws.Range("S6").Value = 0.4140625
$ws.Range("B7").Value = 0.09865470852017937
$ws.Range("D7").Value = 0.0179372197309417
$ws.Range("F7").Value = 0.07174887892376682
$ws.Range("J7").Value = 0.1659192825112108
$ws.Range("O7").Value = 0.02242152466367713
$ws.Range("Q7").Value = 0.1390134529147982
$ws.Range("R7").Value = 0.1255605381165919
$ws.Range("S7").Value = 0.3587443946188341
$ws.Range("B8").Value = 0.07829977628635347
$ws.Range("D8").Value = 0.008948545861297539
$ws.Range("F8").Value = 0.07158836689038031
$ws.Range("J8").Value = 0.0894854586129754
$ws.Range("O8").Value = 0.01789709172259508
$ws.Range("Q8").Value = 0.1722595078299776
$ws.Range("R8").Value = 0.1521252796420582
$ws.Range("S8").Value = 0.4093959731543624
$ws.Range("B9").Value = 0.1
$ws.Range("D9").Value = 0.005
$ws.Range("F9").Value = 0.08
$ws.Range("J9").Value = 0.105
$ws.Range("O9").Value = 0.025
$ws.Range("Q9").Value = 0.14
$ws.Range("R9").Value = 0.175
$ws.Range("S9").Value = 0.37
$ws.Range("B10").Value = 0.1069182389937107
$ws.Range("D10").Value = 0.01437556154537287
$ws.Range("F10").Value = 0.09613656783468104
$ws.Range("J10").Value = 0.08086253369272237
$ws.Range("O10").Value = 0.01617250673854448
$ws.Range("Q10").Value = 0.2012578616352201
$ws.Range("R10").Value = 0.0862533692722372
$ws.Range("S10").Value = 0.3980233602875112
$ws.Range("G11").Value = 0.1358024691358025
$ws.Range("J11").Value = 0.08641975308641975
$ws.Range("K11").Value = 0.1851851851851852
$ws.Range("L11").Value = 0.5833333333333334
$ws.Range("S11").Value = 0.009259259259259259
$ws.Range("G12").Value = 0.7317073170731707
$ws.Range("J12").Value = 0.1707317073170732
$ws.Range("K12").Value = 0.00975609756097561
$ws.Range("L12").Value = 0.04878048780487805
$ws.Range("S12").Value = 0.03902439024390244
$ws.Range("G13").Value = 0.7872340425531915
$ws.Range("J13").Value = 0.1702127659574468
$ws.Range("S13").Value = 0.0425531914893617
$ws.Range("F15").Value = 0.0131578947368421
$ws.Range("H15").Value = 0.2105263157894737
$ws.Range("I15").Value = 0.06578947368421052
$ws.Range("J15").Value = 0.3421052631578947
$ws.Range("K15").Value = 0.05263157894736842
$ws.Range("M15").Value = 0.004385964912280702
$ws.Range("N15").Value = 0.008771929824561403
$ws.Range("O15").Value = 0.05263157894736842
$ws.Range("S15").Value = 0.25
$ws.Range("F16").Value = 0.03773584905660377
$ws.Range("H16").Value = 0.1635220125786163
$ws.Range("I16").Value = 0.08176100628930817
$ws.Range("J16").Value = 0.4025157232704403
$ws.Range("K16").Value = 0.1132075471698113
$ws.Range("M16").Value = 0.02515723270440252
$ws.Range("N16").Value = 0.006289308176100629
$ws.Range("O16").Value = 0.06289308176100629
$ws.Range("S16").Value = 0.1069182389937107
$ws.Range("F17").Value = 0.01566579634464752
$ws.Range("H17").Value = 0.2088772845953002
$ws.Range("I17").Value = 0.08355091383812011
$ws.Range("J17").Value = 0.3864229765013055
$ws.Range("K17").Value = 0.1253263707571801
$ws.Range("M17").Value = 0.01566579634464752
$ws.Range("O17").Value = 0.04960835509138381
$ws.Range("S17").Value = 0.1148825065274151
$ws.Range("F18").Value = 0.02390438247011952
$ws.Range("H18").Value = 0.2111553784860558
$ws.Range("I18").Value = 0.08366533864541832
$ws.Range("J18").Value = 0.3386454183266932
$ws.Range("K18").Value = 0.09561752988047809
$ws.Range("M18").Value = 0.02788844621513944
$ws.Range("O18").Value = 0.08764940239043825
$ws.Range("S18").Value = 0.1314741035856574
$ws.Range("F19").Value = 0.02066772655007949
$ws.Range("H19").Value = 0.1987281399046105
$ws.Range("I19").Value = 0.0945945945945946
$ws.Range("J19").Value = 0.3426073131955485
$ws.Range("K19").Value = 0.1208267090620032
$ws.Range("M19").Value = 0.02305246422893482
$ws.Range("N19").Value = 0.000794912559618442
$ws.Range("O19").Value = 0.07710651828298887
$ws.Range("S19").Value = 0.1216216216216216
